$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.063.78"
$ws.Range("E2").Value = "  +4.04%  "
$ws.Range("D3").Value = "2.655.01"
$ws.Range("E3").Value = "  +6.32%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "113.96"
$ws.Range("E5").Value = "  +7.63%  "
$ws.Range("D6").Value = "326.31"
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("D7").Value = "0.529"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("D10").Value = "41.20"
$ws.Range("E10").Value = "  +6.06%  "
$ws.Range("D11").Value = "20.14"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "0.0826"
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "7.40"
$ws.Range("E14").Value = "  +4.22%  "
$ws.Range("D15").Value = "3.072.67"
$ws.Range("E15").Value = "  +6.37%  "
$ws.Range("D16").Value = "2.657.42"
$ws.Range("E16").Value = "  +6.42%  "
$ws.Range("E17").Value = "  +5.65%  "
$ws.Range("D18").Value = "49.999.77"
$ws.Range("E18").Value = "  +4.18%  "
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  +3.38%  "
$ws.Range("D20").Value = "6.78"
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +3.29%  "
$ws.Range("D23").Value = "72.54"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("D24").Value = "276.28"
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").Value = "26.93"
$ws.Range("E26").Value = "  +4.51%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "36.87"
$ws.Range("E29").Value = "  +6.52%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").Value = "50.21"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("D33").Value = "5.50"
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("D34").Value = "0.0817"
$ws.Range("E34").Value = "  +5.70%  "
$ws.Range("D35").Value = "19.49"
$ws.Range("E35").Value = "  +2.02%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "5.04"
$ws.Range("E37").Value = "  +9.45%  "
$ws.Range("D38").Value = "2.08"
$ws.Range("E38").Value = "  +6.78%  "
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  +8.71%  "
$ws.Range("D40").Value = "124.81"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "22.22"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "2.22"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "0.0319"
$ws.Range("E44").Value = "  +5.83%  "
$ws.Range("D45").Value = "2.087.06"
$ws.Range("E45").Value = "  +4.29%  "
$ws.Range("D46").Value = "3.35"
$ws.Range("E46").Value = "  +6.15%  "
$ws.Range("D47").Value = "2.26"
$ws.Range("E47").Value = "  +13.42%  "
$ws.Range("E48").Value = "  +4.58%  "
$ws.Range("D49").Value = "9.13"
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "5.35"
$ws.Range("E50").Value = "  +2.97%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "60.55"
$ws.Range("E51").Value = "  +6.91%  "
